$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 99275.06305729281
$ws.Range("C2").Value = 103576.4522698135
$ws.Range("D2").Value = 142715.8049429919
